# Apply updated policy-effectiveness values to the "Entertainment(Indoor)" (F)
# and "Industries" (H) columns, mirroring the change captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Entertainment(Indoor)") lag weighting updates.
$ws.Range("F16:F24").Value   = 0.1428571428571428
$ws.Range("F25:F129").Value  = 0.8571428571428571
$ws.Range("F130:F136").Value = 0.5714285714285714

# Column H ("Industries") values reset to 0 for rows 25-101.
$ws.Range("H25:H101").Value = 0

Write-Output "Updated F16:F136 and H25:H101"
